# Scheduled runner update: refresh cached market-board figures
# (currentAveragePrice / LevePrice* / LeveProfit* columns) per-sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 87
$ws.Cells.Item(87, 8).Value = 11541.404
$ws.Cells.Item(87, 10).Value = 11541.404
$ws.Cells.Item(87, 12).Value = 11541.404
$ws.Cells.Item(87, 14).Value = -14037.404
# Row 90
$ws.Cells.Item(90, 8).Value = 11541.404
$ws.Cells.Item(90, 10).Value = 11541.404
$ws.Cells.Item(90, 12).Value = 34624.212
$ws.Cells.Item(90, 14).Value = -47104.212
# Row 98
$ws.Cells.Item(98, 8).Value = 6569637.5
$ws.Cells.Item(98, 9).Value = 9031509
$ws.Cells.Item(98, 10).Value = 4646.5
$ws.Cells.Item(98, 11).Value = 9031509
$ws.Cells.Item(98, 12).Value = 4646.5
$ws.Cells.Item(98, 13).Value = -9030011
$ws.Cells.Item(98, 14).Value = -7642.5
# Row 122
$ws.Cells.Item(122, 8).Value = 6569637.5
$ws.Cells.Item(122, 9).Value = 9031509
$ws.Cells.Item(122, 10).Value = 4646.5
$ws.Cells.Item(122, 11).Value = 27094527
$ws.Cells.Item(122, 12).Value = 13939.5
$ws.Cells.Item(122, 13).Value = -27092077
$ws.Cells.Item(122, 14).Value = -18839.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 346308.7
$ws.Cells.Item(32, 9).Value = 3406.2463
$ws.Cells.Item(32, 10).Value = 1923659.9
$ws.Cells.Item(32, 11).Value = 3406.2463
$ws.Cells.Item(32, 12).Value = 1923659.9
$ws.Cells.Item(32, 13).Value = -3119.2463
$ws.Cells.Item(32, 14).Value = -1924233.9
# Row 36
$ws.Cells.Item(36, 8).Value = 6880.8
$ws.Cells.Item(36, 9).Value = 6880.8
$ws.Cells.Item(36, 11).Value = 6880.8
$ws.Cells.Item(36, 13).Value = -6534.8
# Row 42
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = $null
$ws.Cells.Item(42, 14).Value = 0
# Row 61
$ws.Cells.Item(61, 8).Value = 3215.5833
$ws.Cells.Item(61, 9).Value = 2780.647
$ws.Cells.Item(61, 10).Value = 4271.857
$ws.Cells.Item(61, 11).Value = 2780.647
$ws.Cells.Item(61, 12).Value = 4271.857
$ws.Cells.Item(61, 13).Value = -2568.647
$ws.Cells.Item(61, 14).Value = -4695.857
# Row 74
$ws.Cells.Item(74, 8).Value = 1074
$ws.Cells.Item(74, 9).Value = 962.3077
$ws.Cells.Item(74, 10).Value = 1281.4286
$ws.Cells.Item(74, 11).Value = 962.3077
$ws.Cells.Item(74, 12).Value = 1281.4286
$ws.Cells.Item(74, 13).Value = -88.30769999999995
$ws.Cells.Item(74, 14).Value = -3029.4286
# Row 77
$ws.Cells.Item(77, 8).Value = 1074
$ws.Cells.Item(77, 9).Value = 962.3077
$ws.Cells.Item(77, 10).Value = 1281.4286
$ws.Cells.Item(77, 11).Value = 4811.5385
$ws.Cells.Item(77, 12).Value = 6407.143
$ws.Cells.Item(77, 13).Value = -443.5384999999997
$ws.Cells.Item(77, 14).Value = -15143.143
# Row 136
$ws.Cells.Item(136, 8).Value = 3215.5833
$ws.Cells.Item(136, 9).Value = 2780.647
$ws.Cells.Item(136, 10).Value = 4271.857
$ws.Cells.Item(136, 11).Value = 8341.940999999999
$ws.Cells.Item(136, 12).Value = 12815.571
$ws.Cells.Item(136, 13).Value = -5791.940999999999
$ws.Cells.Item(136, 14).Value = -17915.571

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 9709.071
$ws.Cells.Item(134, 9).Value = 3978.0952
$ws.Cells.Item(134, 10).Value = 26902
$ws.Cells.Item(134, 11).Value = 11934.2856
$ws.Cells.Item(134, 12).Value = 80706
$ws.Cells.Item(134, 13).Value = -9399.285600000001
$ws.Cells.Item(134, 14).Value = -85776

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Cells.Item(20, 8).Value = 49172.086
$ws.Cells.Item(20, 10).Value = 49172.086
$ws.Cells.Item(20, 12).Value = 49172.086
$ws.Cells.Item(20, 14).Value = -49644.086
# Row 30
$ws.Cells.Item(30, 8).Value = 49172.086
$ws.Cells.Item(30, 10).Value = 49172.086
$ws.Cells.Item(30, 12).Value = 49172.086
$ws.Cells.Item(30, 14).Value = -49354.086
# Row 94
$ws.Cells.Item(94, 8).Value = 22733668
$ws.Cells.Item(94, 9).Value = 62500530
$ws.Cells.Item(94, 10).Value = 9748
$ws.Cells.Item(94, 11).Value = 62500530
$ws.Cells.Item(94, 12).Value = 9748
$ws.Cells.Item(94, 13).Value = -62500079
$ws.Cells.Item(94, 14).Value = -10650
# Row 128
$ws.Cells.Item(128, 8).Value = 49172.086
$ws.Cells.Item(128, 10).Value = 49172.086
$ws.Cells.Item(128, 12).Value = 49172.086
$ws.Cells.Item(128, 14).Value = -59132.086

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Cells.Item(12, 8).Value = 164.11765
$ws.Cells.Item(12, 9).Value = 1.5
$ws.Cells.Item(12, 10).Value = 554.4
$ws.Cells.Item(12, 11).Value = 4.5
$ws.Cells.Item(12, 12).Value = 1663.2
$ws.Cells.Item(12, 13).Value = 168.5
$ws.Cells.Item(12, 14).Value = -2009.2
# Row 113
$ws.Cells.Item(113, 8).Value = 827.8
$ws.Cells.Item(113, 9).Value = 578.58826
$ws.Cells.Item(113, 11).Value = 1735.76478
$ws.Cells.Item(113, 13).Value = 434.23522
# Row 131
$ws.Cells.Item(131, 8).Value = 7043290
$ws.Cells.Item(131, 9).Value = 957.875
$ws.Cells.Item(131, 10).Value = 7937554.5
$ws.Cells.Item(131, 11).Value = 2873.625
$ws.Cells.Item(131, 12).Value = 23812663.5
$ws.Cells.Item(131, 13).Value = 2166.375
$ws.Cells.Item(131, 14).Value = -23822743.5
# Row 138
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 13).Value = $null

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 58720.55
$ws.Cells.Item(80, 9).Value = 161258.58
$ws.Cells.Item(80, 10).Value = 3507.7693
$ws.Cells.Item(80, 11).Value = 161258.58
$ws.Cells.Item(80, 12).Value = 3507.7693
$ws.Cells.Item(80, 13).Value = -160260.58
$ws.Cells.Item(80, 14).Value = -5503.7693
# Row 83
$ws.Cells.Item(83, 8).Value = 58720.55
$ws.Cells.Item(83, 9).Value = 161258.58
$ws.Cells.Item(83, 10).Value = 3507.7693
$ws.Cells.Item(83, 11).Value = 806292.8999999999
$ws.Cells.Item(83, 12).Value = 17538.8465
$ws.Cells.Item(83, 13).Value = -801300.8999999999
$ws.Cells.Item(83, 14).Value = -27522.8465

$ws = $wb.Worksheets.Item("LTW")
# Row 80
$ws.Cells.Item(80, 8).Value = 20888.889
$ws.Cells.Item(80, 10).Value = 20888.889
$ws.Cells.Item(80, 12).Value = 20888.889
$ws.Cells.Item(80, 14).Value = -23134.889
# Row 83
$ws.Cells.Item(83, 8).Value = 20888.889
$ws.Cells.Item(83, 10).Value = 20888.889
$ws.Cells.Item(83, 12).Value = 62666.667
$ws.Cells.Item(83, 14).Value = -73898.667
# Row 132
$ws.Cells.Item(132, 8).Value = 3347541.2
$ws.Cells.Item(132, 9).Value = 6411837.5
$ws.Cells.Item(132, 10).Value = 4672.5454
$ws.Cells.Item(132, 11).Value = 19235512.5
$ws.Cells.Item(132, 12).Value = 14017.6362
$ws.Cells.Item(132, 13).Value = -19232982.5
$ws.Cells.Item(132, 14).Value = -19077.6362

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81, 8).Value = 3977.353
$ws.Cells.Item(81, 9).Value = 3484.1177
$ws.Cells.Item(81, 10).Value = 4470.5884
$ws.Cells.Item(81, 11).Value = 6968.2354
$ws.Cells.Item(81, 12).Value = 8941.176799999999
$ws.Cells.Item(81, 13).Value = -5907.2354
$ws.Cells.Item(81, 14).Value = -11063.1768
# Row 84
$ws.Cells.Item(84, 8).Value = 3977.353
$ws.Cells.Item(84, 9).Value = 3484.1177
$ws.Cells.Item(84, 10).Value = 4470.5884
$ws.Cells.Item(84, 11).Value = 34841.177
$ws.Cells.Item(84, 12).Value = 44705.884
$ws.Cells.Item(84, 13).Value = -29537.177
$ws.Cells.Item(84, 14).Value = -55313.884
# Row 136
$ws.Cells.Item(136, 8).Value = 77663.234
$ws.Cells.Item(136, 9).Value = 125626.5
$ws.Cells.Item(136, 10).Value = 922
$ws.Cells.Item(136, 11).Value = 376879.5
$ws.Cells.Item(136, 12).Value = 2766
$ws.Cells.Item(136, 13).Value = -374329.5
$ws.Cells.Item(136, 14).Value = -7866
